$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain decimal number need their
# number format forced to Text first, otherwise Excel auto-converts the
# string into a numeric value (losing the intended text representation).
$textCells = @("D5","D6","D11","D13","D14","D18","D20","D21","D23","D24","D31","D35","D37","D38","D41","D43","D45","D46","D50")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "65.869.61"
$ws.Range("E2").Value = "  +1.48%  "
$ws.Range("D3").Value = "2.696.94"
$ws.Range("E3").Value = "  +2.18%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "608.86"
$ws.Range("E5").Value = "  +2.26%  "
$ws.Range("D6").Value = "158.03"
$ws.Range("E6").Value = "  +1.75%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  -0.66%  "
$ws.Range("E9").Value = "  +6.21%  "
$ws.Range("E10").Value = "  +4.08%  "
$ws.Range("D11").Value = "0.402"
$ws.Range("E11").Value = "  +1.19%  "
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("D13").Value = "30.38"
$ws.Range("E13").Value = "  +4.51%  "
$ws.Range("D14").Value = "0.0000201"
$ws.Range("E14").Value = "  +7.95%  "
$ws.Range("D15").Value = "3.180.83"
$ws.Range("E15").Value = "  +2.02%  "
$ws.Range("D16").Value = "65.721.37"
$ws.Range("E16").Value = "  +1.38%  "
$ws.Range("D17").Value = "2.705.80"
$ws.Range("E17").Value = "  +2.40%  "
$ws.Range("D18").Value = "12.65"
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("E19").Value = "  +2.15%  "
$ws.Range("D20").Value = "359.49"
$ws.Range("E20").Value = "  +2.36%  "
$ws.Range("D21").Value = "7.56"
$ws.Range("E21").Value = "  +3.75%  "
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("D23").Value = "70.72"
$ws.Range("E23").Value = "  +4.20%  "
$ws.Range("D24").Value = "9.85"
$ws.Range("E24").Value = "  +3.66%  "
$ws.Range("E25").Value = "  +14.09%  "
$ws.Range("E26").Value = "  -1.33%  "
$ws.Range("E27").Value = "  +2.90%  "
$ws.Range("E28").Value = "  +5.00%  "
$ws.Range("E29").Value = "  +3.56%  "
$ws.Range("E30").Value = "  +5.29%  "
$ws.Range("D31").Value = "541.22"
$ws.Range("E31").Value = "  +5.87%  "
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("E33").Value = "  +1.88%  "
$ws.Range("E34").Value = "  +6.00%  "
$ws.Range("D35").Value = "5.39"
$ws.Range("E35").Value = "  -3.49%  "
$ws.Range("E36").Value = "  +1.90%  "
$ws.Range("D37").Value = "20.84"
$ws.Range("E37").Value = "  +3.28%  "
$ws.Range("D38").Value = "163.41"
$ws.Range("E38").Value = "  -0.75%  "
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("D41").Value = "172.44"
$ws.Range("E41").Value = "  +4.65%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").Value = "42.62"
$ws.Range("E43").Value = "  +0.94%  "
$ws.Range("E44").Value = "  +2.49%  "
$ws.Range("D45").Value = "0.0614"
$ws.Range("E45").Value = "  +0.40%  "
$ws.Range("D46").Value = "23.61"
$ws.Range("E46").Value = "  +3.18%  "
$ws.Range("E47").Value = "  +4.40%  "
$ws.Range("E49").Value = "  +1.68%  "
$ws.Range("D50").Value = "21.09"
$ws.Range("E50").Value = "  +9.36%  "
$ws.Range("E51").Value = "  +1.34%  "
